$wb = $excel.ActiveWorkbook

# --- "Current" sheet: move selection back to the top (A2) ---
$wsCurrent = $wb.Worksheets.Item("Current")
$wsCurrent.Range("A2").Select()

# --- "v2.0+" sheet: insert a new bullet row about e-mail identifiers ---
$ws = $wb.Worksheets.Item("v2.0+")
$ws.Rows(55).Insert()
$ws.Range("A55").Value = "o`u{00A0}`u{00A0} STRETCH `u{2013} change results to not use e-mails as identifiers."
$ws.Rows(55).RowHeight = 45

# Update the view/selection for this sheet as well
$ws.Range("B55").Select()
